$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain text (e.g. "0.4680", "1.842.38") rather than
# real numbers. Several of the new values look numeric, so Excel would silently
# coerce them (dropping trailing zeros, switching to scientific notation, etc.)
# if assigned directly. Mark those cells as Text first so the literal string survives.
$textCells = @("D5","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Cell value updates (latest crypto snapshot) ---
$ws.Range("D2").Value = '26.850.31'
$ws.Range("E2").Value = '  +1.60%  '
$ws.Range("D3").Value = '1.838.11'
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").Value = '309.71'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '0.4680'
$ws.Range("E7").Value = '  +3.43%  '
$ws.Range("D8").Value = '0.3624'
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("D9").Value = '0.07170'
$ws.Range("E9").Value = '  +1.46%  '
$ws.Range("D10").Value = '0.9363'
$ws.Range("E10").Value = '  +5.37%  '
$ws.Range("D11").Value = '19.57'
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D12").Value = '0.07680'
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("D13").Value = '1.842.38'
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("D14").Value = '5.283'
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("D15").Value = '6.376'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '88.04'
$ws.Range("E16").Value = '  +3.28%  '
$ws.Range("D17").Value = '1.010'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").Value = '0.000008572'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").Value = '26.851.65'
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").Value = '14.32'
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("D22").Value = '5.031'
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("D23").Value = '10.62'
$ws.Range("E23").Value = '  +1.04%  '
$ws.Range("D24").Value = '1.917'
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("D25").Value = '152.31'
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("D26").Value = '18.00'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").Value = '2.006'
$ws.Range("E27").Value = '  -1.84%  '
$ws.Range("D28").Value = '113.95'
$ws.Range("E28").Value = '  +1.78%  '
$ws.Range("D29").Value = '4.909'
$ws.Range("E29").Value = '  +1.13%  '
$ws.Range("D30").Value = '0.08843'
$ws.Range("E30").Value = '  +1.80%  '
$ws.Range("D31").Value = '3.161'
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("D32").Value = '2.856'
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").Value = '1.183'
$ws.Range("E33").Value = '  +7.35%  '
$ws.Range("D34").Value = '0.7456'
$ws.Range("E34").Value = '  +3.78%  '
$ws.Range("D35").Value = '4.458'
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("D36").Value = '1.088'
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").Value = '2.988'
$ws.Range("E37").Value = '  +3.23%  '
$ws.Range("D38").Value = '0.01931'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = '0.05143'
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5108'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.910'
$ws.Range("E41").Value = '  +1.76%  '
$ws.Range("D42").Value = '0.1515'
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").Value = '8.148'
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("D44").Value = '0.4704'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("E45").Value = '  +2.59%  '
$ws.Range("D46").Value = '1.009'
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("D47").Value = '99.69'
$ws.Range("E47").Value = '  -0.66%  '
$ws.Range("D48").Value = '1.588'
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("D50").Value = '63.97'
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("D51").Value = '36.11'
$ws.Range("E51").Value = '  +0.03%  '
